$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows directly below the header row (row 1), pushing the
# existing data rows down by 9 (old row 2 -> new row 11, ... old row 21 -> new row 30)
$ws.Rows("2:10").Insert()

# Insert() copies the formatting of the row above (the bold header) onto the
# freshly inserted rows; clear that so the new data rows stay unstyled like the rest
$ws.Rows("2:10").ClearFormats()

# Populate the newly inserted rows with the new accelerometer readings
$ws.Cells.Item(2, 1).Value = -1.231432914733887
$ws.Cells.Item(2, 2).Value = 2.921578645706177
$ws.Cells.Item(2, 3).Value = 1.135899901390076

$ws.Cells.Item(3, 1).Value = -1.530599117279053
$ws.Cells.Item(3, 2).Value = 3.089309453964233
$ws.Cells.Item(3, 3).Value = 1.274296164512634

$ws.Cells.Item(4, 1).Value = -1.078460693359375
$ws.Cells.Item(4, 2).Value = 3.193733692169189
$ws.Cells.Item(4, 3).Value = 1.281612634658814

$ws.Cells.Item(5, 1).Value = -0.9324893951416016
$ws.Cells.Item(5, 2).Value = 2.998547315597534
$ws.Cells.Item(5, 3).Value = 0.8108012080192566

$ws.Cells.Item(6, 1).Value = -1.057272911071777
$ws.Cells.Item(6, 2).Value = 2.989111185073853
$ws.Cells.Item(6, 3).Value = 0.8697453737258911

$ws.Cells.Item(7, 1).Value = -1.247255325317383
$ws.Cells.Item(7, 2).Value = 3.032690763473511
$ws.Cells.Item(7, 3).Value = 0.9775734543800354

$ws.Cells.Item(8, 1).Value = -1.120566368103027
$ws.Cells.Item(8, 2).Value = 3.040028095245361
$ws.Cells.Item(8, 3).Value = 0.9562293887138368

$ws.Cells.Item(9, 1).Value = -1.322433471679688
$ws.Cells.Item(9, 2).Value = 3.130712985992432
$ws.Cells.Item(9, 3).Value = 1.133776545524597

$ws.Cells.Item(10, 1).Value = -1.53396463394165
$ws.Cells.Item(10, 2).Value = 3.08948278427124
$ws.Cells.Item(10, 3).Value = 1.223363161087036

# Append one additional new row of data at the end of the table
$ws.Cells.Item(31, 1).Value = -0.5697603225708008
$ws.Cells.Item(31, 2).Value = 2.878552436828613
$ws.Cells.Item(31, 3).Value = 1.078300476074219
